# Auto-generated edit script: apply scheduled market-data refresh to Sheets/Rafflesia_Profits.xlsx
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) across all 8 sheets
# to reflect the latest pulled market data. Cells that should be blank (no applicable price)
# are cleared rather than zeroed, matching upstream behavior.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1029.5  # H17: 953 -> 1029.5
$ws.Cells.Item(17, 10).Value = 1029.5  # J17: 953 -> 1029.5
$ws.Cells.Item(17, 12).Value = 3088.5  # L17: 2859 -> 3088.5
$ws.Cells.Item(17, 14).Value = -3424.5  # N17: -3195 -> -3424.5

$ws.Cells.Item(26, 8).Value = 18015  # H26: 15347.667 -> 18015
$ws.Cells.Item(26, 9).Value = 0  # I26: 10013 -> 0
$ws.Cells.Item(26, 11).Value = 0  # K26: 10013 -> 0
$ws.Cells.Item(26, 13).ClearContents()  # M26: -9669 -> (blank)

$ws.Cells.Item(33, 8).Value = 549.41174  # H33: 519.44446 -> 549.41174
$ws.Cells.Item(33, 9).Value = 156  # I33: 146.875 -> 156
$ws.Cells.Item(33, 11).Value = 156  # K33: 146.875 -> 156
$ws.Cells.Item(33, 13).Value = 73  # M33: 82.125 -> 73

$ws.Cells.Item(40, 8).Value = 0  # H40: 2749.5 -> 0
$ws.Cells.Item(40, 9).Value = 0  # I40: 1499 -> 0
$ws.Cells.Item(40, 10).Value = 0  # J40: 4000 -> 0
$ws.Cells.Item(40, 11).Value = 0  # K40: 1499 -> 0
$ws.Cells.Item(40, 12).Value = 0  # L40: 4000 -> 0
$ws.Cells.Item(40, 13).ClearContents()  # M40: -1324 -> (blank)
$ws.Cells.Item(40, 14).ClearContents()  # N40: -4350 -> (blank)

$ws.Cells.Item(55, 8).Value = 201.8  # H55: 230.66667 -> 201.8
$ws.Cells.Item(55, 9).Value = 125  # I55: 200 -> 125
$ws.Cells.Item(55, 10).Value = 221  # J55: 236.8 -> 221
$ws.Cells.Item(55, 11).Value = 125  # K55: 200 -> 125
$ws.Cells.Item(55, 12).Value = 221  # L55: 236.8 -> 221
$ws.Cells.Item(55, 13).Value = 89  # M55: 14 -> 89
$ws.Cells.Item(55, 14).Value = -649  # N55: -664.8 -> -649

$ws.Cells.Item(80, 8).Value = 1316.6666  # H80: 2000 -> 1316.6666
$ws.Cells.Item(80, 10).Value = 975  # J80: 0 -> 975
$ws.Cells.Item(80, 12).Value = 2925  # L80: 0 -> 2925
$ws.Cells.Item(80, 14).Value = -4921  # N80: None -> -4921

$ws.Cells.Item(83, 8).Value = 1316.6666  # H83: 2000 -> 1316.6666
$ws.Cells.Item(83, 10).Value = 975  # J83: 0 -> 975
$ws.Cells.Item(83, 12).Value = 8775  # L83: 0 -> 8775
$ws.Cells.Item(83, 14).Value = -18759  # N83: None -> -18759

$ws.Cells.Item(98, 8).Value = 83334490  # H98: 100001290 -> 83334490
$ws.Cells.Item(98, 9).Value = 100001150  # I98: 125001310 -> 100001150
$ws.Cells.Item(98, 11).Value = 100001150  # K98: 125001310 -> 100001150
$ws.Cells.Item(98, 13).Value = -99999652  # M98: -124999812 -> -99999652

$ws.Cells.Item(99, 8).Value = 1150  # H99: 90 -> 1150
$ws.Cells.Item(99, 9).Value = 0  # I99: 90 -> 0
$ws.Cells.Item(99, 10).Value = 1150  # J99: 0 -> 1150
$ws.Cells.Item(99, 11).Value = 0  # K99: 270 -> 0
$ws.Cells.Item(99, 12).Value = 3450  # L99: 0 -> 3450
$ws.Cells.Item(99, 13).ClearContents()  # M99: 1228 -> (blank)
$ws.Cells.Item(99, 14).Value = -6446  # N99: None -> -6446

$ws.Cells.Item(101, 8).Value = 4000  # H101: 0 -> 4000
$ws.Cells.Item(101, 9).Value = 4000  # I101: 0 -> 4000
$ws.Cells.Item(101, 11).Value = 12000  # K101: 0 -> 12000
$ws.Cells.Item(101, 13).Value = -10378  # M101: None -> -10378

$ws.Cells.Item(122, 8).Value = 83334490  # H122: 100001290 -> 83334490
$ws.Cells.Item(122, 9).Value = 100001150  # I122: 125001310 -> 100001150
$ws.Cells.Item(122, 11).Value = 300003450  # K122: 375003930 -> 300003450
$ws.Cells.Item(122, 13).Value = -300001000  # M122: -375001480 -> -300001000

$ws.Cells.Item(137, 8).Value = 1725.6  # H137: 1784.8948 -> 1725.6
$ws.Cells.Item(137, 9).Value = 1383.1765  # I137: 1432.1875 -> 1383.1765
$ws.Cells.Item(137, 11).Value = 4149.529500000001  # K137: 4296.5625 -> 4149.529500000001
$ws.Cells.Item(137, 13).Value = -1599.529500000001  # M137: -1746.5625 -> -1599.529500000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1088.7368  # H32: 1109.7894 -> 1088.7368
$ws.Cells.Item(32, 9).Value = 771  # I32: 796 -> 771
$ws.Cells.Item(32, 11).Value = 771  # K32: 796 -> 771
$ws.Cells.Item(32, 13).Value = -484  # M32: -509 -> -484

$ws.Cells.Item(97, 8).Value = 1117  # H97: 1271.25 -> 1117
$ws.Cells.Item(97, 10).Value = 638.5  # J97: 777 -> 638.5
$ws.Cells.Item(97, 12).Value = 638.5  # L97: 777 -> 638.5
$ws.Cells.Item(97, 14).Value = -1630.5  # N97: -1769 -> -1630.5

$ws.Cells.Item(101, 8).Value = 25000  # H101: 22666.666 -> 25000
$ws.Cells.Item(101, 10).Value = 25000  # J101: 22666.666 -> 25000
$ws.Cells.Item(101, 12).Value = 25000  # L101: 22666.666 -> 25000
$ws.Cells.Item(101, 14).Value = -31490  # N101: -29156.666 -> -31490

$ws.Cells.Item(137, 8).Value = 20000  # H137: 0 -> 20000
$ws.Cells.Item(137, 9).Value = 20000  # I137: 0 -> 20000
$ws.Cells.Item(137, 11).Value = 20000  # K137: 0 -> 20000
$ws.Cells.Item(137, 13).Value = -14900  # M137: None -> -14900

$ws.Cells.Item(139, 8).Value = 89904.336  # H139: 105000 -> 89904.336
$ws.Cells.Item(139, 10).Value = 89904.336  # J139: 105000 -> 89904.336
$ws.Cells.Item(139, 12).Value = 89904.336  # L139: 105000 -> 89904.336
$ws.Cells.Item(139, 14).Value = -100184.336  # N139: -115280 -> -100184.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(8, 8).Value = 200  # H8: 0 -> 200
$ws.Cells.Item(8, 9).Value = 200  # I8: 0 -> 200
$ws.Cells.Item(8, 11).Value = 200  # K8: 0 -> 200
$ws.Cells.Item(8, 13).Value = -60  # M8: None -> -60

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 328  # H22: 240.4 -> 328
$ws.Cells.Item(22, 9).Value = 299.6  # I22: 199.5 -> 299.6
$ws.Cells.Item(22, 10).Value = 399  # J22: 404 -> 399
$ws.Cells.Item(22, 11).Value = 299.6  # K22: 199.5 -> 299.6
$ws.Cells.Item(22, 12).Value = 399  # L22: 404 -> 399
$ws.Cells.Item(22, 13).Value = 50.39999999999998  # M22: 150.5 -> 50.39999999999998
$ws.Cells.Item(22, 14).Value = -1099  # N22: -1104 -> -1099

$ws.Cells.Item(41, 8).Value = 31011.54  # H41: 28939.285 -> 31011.54
$ws.Cells.Item(41, 9).Value = 2500  # I41: 2250 -> 2500
$ws.Cells.Item(41, 11).Value = 2500  # K41: 2250 -> 2500
$ws.Cells.Item(41, 13).Value = -2072  # M41: -1822 -> -2072

$ws.Cells.Item(59, 8).Value = 0  # H59: 120000 -> 0
$ws.Cells.Item(59, 10).Value = 0  # J59: 120000 -> 0
$ws.Cells.Item(59, 12).Value = 0  # L59: 120000 -> 0
$ws.Cells.Item(59, 14).ClearContents()  # N59: -122290 -> (blank)

$ws.Cells.Item(132, 8).Value = 5286.5  # H132: 5699.136 -> 5286.5
$ws.Cells.Item(132, 9).Value = 2005.9231  # I132: 2234.7273 -> 2005.9231
$ws.Cells.Item(132, 11).Value = 6017.7693  # K132: 6704.1819 -> 6017.7693
$ws.Cells.Item(132, 13).Value = -3487.7693  # M132: -4174.1819 -> -3487.7693

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 215000  # H4: 250000 -> 215000
$ws.Cells.Item(4, 10).Value = 215000  # J4: 250000 -> 215000
$ws.Cells.Item(4, 12).Value = 645000  # L4: 750000 -> 645000
$ws.Cells.Item(4, 14).Value = -645224  # N4: -750224 -> -645224

$ws.Cells.Item(23, 8).Value = 526.1667  # H23: 520.7143 -> 526.1667
$ws.Cells.Item(23, 10).Value = 556.6  # J23: 545.1667 -> 556.6
$ws.Cells.Item(23, 12).Value = 1669.8  # L23: 1635.5001 -> 1669.8
$ws.Cells.Item(23, 14).Value = -2139.8  # N23: -2105.5001 -> -2139.8

$ws.Cells.Item(52, 8).Value = 1500  # H52: 0 -> 1500
$ws.Cells.Item(52, 10).Value = 1500  # J52: 0 -> 1500
$ws.Cells.Item(52, 12).Value = 4500  # L52: 0 -> 4500
$ws.Cells.Item(52, 14).Value = -5032  # N52: None -> -5032

$ws.Cells.Item(115, 8).Value = 1999.6666  # H115: 0 -> 1999.6666
$ws.Cells.Item(115, 9).Value = 2000  # I115: 0 -> 2000
$ws.Cells.Item(115, 10).Value = 1999  # J115: 0 -> 1999
$ws.Cells.Item(115, 11).Value = 6000  # K115: 0 -> 6000
$ws.Cells.Item(115, 12).Value = 5997  # L115: 0 -> 5997
$ws.Cells.Item(115, 13).Value = -4825  # M115: None -> -4825
$ws.Cells.Item(115, 14).Value = -8347  # N115: None -> -8347

$ws.Cells.Item(131, 8).Value = 1665.5264  # H131: 1678.85 -> 1665.5264
$ws.Cells.Item(131, 10).Value = 1921.8  # J131: 1922.4375 -> 1921.8
$ws.Cells.Item(131, 12).Value = 5765.4  # L131: 5767.3125 -> 5765.4
$ws.Cells.Item(131, 14).Value = -15845.4  # N131: -15847.3125 -> -15845.4

$ws.Cells.Item(140, 8).Value = 1650  # H140: 1775 -> 1650
$ws.Cells.Item(140, 9).Value = 1650  # I140: 1700 -> 1650
$ws.Cells.Item(140, 10).Value = 0  # J140: 2000 -> 0
$ws.Cells.Item(140, 11).Value = 4950  # K140: 5100 -> 4950
$ws.Cells.Item(140, 12).Value = 0  # L140: 6000 -> 0
$ws.Cells.Item(140, 13).Value = 230  # M140: 80 -> 230
$ws.Cells.Item(140, 14).ClearContents()  # N140: -16360 -> (blank)

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 721.25  # H2: 820.7143 -> 721.25
$ws.Cells.Item(2, 9).Value = 928.3333  # I2: 1109 -> 928.3333
$ws.Cells.Item(2, 11).Value = 928.3333  # K2: 1109 -> 928.3333
$ws.Cells.Item(2, 13).Value = -815.3333  # M2: -996 -> -815.3333

$ws.Cells.Item(102, 8).Value = 983.75  # H102: 1588 -> 983.75
$ws.Cells.Item(102, 9).Value = 1048.3334  # I102: 1588 -> 1048.3334
$ws.Cells.Item(102, 10).Value = 790  # J102: 0 -> 790
$ws.Cells.Item(102, 11).Value = 1048.3334  # K102: 1588 -> 1048.3334
$ws.Cells.Item(102, 12).Value = 790  # L102: 0 -> 790
$ws.Cells.Item(102, 13).Value = 573.6666  # M102: 34 -> 573.6666
$ws.Cells.Item(102, 14).Value = -4034  # N102: None -> -4034

$ws.Cells.Item(132, 8).Value = 2322.4  # H132: 3133.3333 -> 2322.4
$ws.Cells.Item(132, 9).Value = 1778  # I132: 2450 -> 1778
$ws.Cells.Item(132, 11).Value = 5334  # K132: 7350 -> 5334
$ws.Cells.Item(132, 13).Value = -2804  # M132: -4820 -> -2804

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 999.5  # H7: 0 -> 999.5
$ws.Cells.Item(7, 9).Value = 999.5  # I7: 0 -> 999.5
$ws.Cells.Item(7, 11).Value = 999.5  # K7: 0 -> 999.5
$ws.Cells.Item(7, 13).Value = -887.5  # M7: None -> -887.5

$ws.Cells.Item(126, 8).Value = 999.5  # H126: 0 -> 999.5
$ws.Cells.Item(126, 9).Value = 999.5  # I126: 0 -> 999.5
$ws.Cells.Item(126, 11).Value = 2998.5  # K126: 0 -> 2998.5
$ws.Cells.Item(126, 13).Value = -528.5  # M126: None -> -528.5

$ws.Cells.Item(132, 8).Value = 8000  # H132: 2000 -> 8000
$ws.Cells.Item(132, 9).Value = 8000  # I132: 2000 -> 8000
$ws.Cells.Item(132, 11).Value = 24000  # K132: 6000 -> 24000
$ws.Cells.Item(132, 13).Value = -21470  # M132: -3470 -> -21470

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(38, 8).Value = 9062.5  # H38: 7541.6665 -> 9062.5
$ws.Cells.Item(38, 9).Value = 9062.5  # I38: 7541.6665 -> 9062.5
$ws.Cells.Item(38, 11).Value = 9062.5  # K38: 7541.6665 -> 9062.5
$ws.Cells.Item(38, 13).Value = -8589.5  # M38: -7068.6665 -> -8589.5

$ws.Cells.Item(81, 8).Value = 800  # H81: 700 -> 800
$ws.Cells.Item(81, 9).Value = 800  # I81: 700 -> 800
$ws.Cells.Item(81, 11).Value = 1600  # K81: 1400 -> 1600
$ws.Cells.Item(81, 13).Value = -539  # M81: -339 -> -539

$ws.Cells.Item(84, 8).Value = 800  # H84: 700 -> 800
$ws.Cells.Item(84, 9).Value = 800  # I84: 700 -> 800
$ws.Cells.Item(84, 11).Value = 8000  # K84: 7000 -> 8000
$ws.Cells.Item(84, 13).Value = -2696  # M84: -1696 -> -2696

$ws.Cells.Item(104, 8).Value = 29000  # H104: 17966.666 -> 29000
$ws.Cells.Item(104, 10).Value = 29000  # J104: 17966.666 -> 29000
$ws.Cells.Item(104, 12).Value = 29000  # L104: 17966.666 -> 29000
$ws.Cells.Item(104, 14).Value = -35988  # N104: -24954.666 -> -35988

$ws.Cells.Item(132, 8).Value = 3759.7144  # H132: 4169.6665 -> 3759.7144
$ws.Cells.Item(132, 9).Value = 3763.3333  # I132: 4256 -> 3763.3333
$ws.Cells.Item(132, 11).Value = 11289.9999  # K132: 12768 -> 11289.9999
$ws.Cells.Item(132, 13).Value = -8759.999899999999  # M132: -10238 -> -8759.999899999999
